# Apply the 'Natmi following Dr Hou advice' update to Sheet1:
# the ligand/receptor data table grows from 8 to 12 data rows (A2:T9 -> A2:T13)
# and every data row's values are recomputed.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$arr = New-Object 'object[,]' 12,20
$arr[0,0] = "ECs"
$arr[0,1] = "Slit2"
$arr[0,2] = "App"
$arr[0,3] = "ECs"
$arr[0,4] = 2
$arr[0,5] = 0.6666666666666666
$arr[0,6] = 0.143896
$arr[0,7] = 0.431688
$arr[0,8] = 0.02807111181859822
$arr[0,9] = 0.02807111181859822
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 110.8604276666667
$arr[0,13] = 332.581283
$arr[0,14] = 0.2509786052589675
$arr[0,15] = 0.2509786052589675
$arr[0,16] = 15.95237209952267
$arr[0,17] = 143.571348895704
$arr[0,18] = 0.0070452484923003
$arr[0,19] = 0.0070452484923003
$arr[1,0] = "ECs"
$arr[1,1] = "Slit2"
$arr[1,2] = "App"
$arr[1,3] = "FAPs"
$arr[1,4] = 2
$arr[1,5] = 0.6666666666666666
$arr[1,6] = 0.143896
$arr[1,7] = 0.431688
$arr[1,8] = 0.02807111181859822
$arr[1,9] = 0.02807111181859822
$arr[1,10] = 3
$arr[1,11] = 1
$arr[1,12] = 184.841802
$arr[1,13] = 554.525406
$arr[1,14] = 0.4184661617850055
$arr[1,15] = 0.4184661617850055
$arr[1,16] = 26.597995940592
$arr[1,17] = 239.381963465328
$arr[1,18] = 0.0117468104197665
$arr[1,19] = 0.0117468104197665
$arr[2,0] = "ECs"
$arr[2,1] = "Slit2"
$arr[2,2] = "App"
$arr[2,3] = "M2"
$arr[2,4] = 2
$arr[2,5] = 0.6666666666666666
$arr[2,6] = 0.143896
$arr[2,7] = 0.431688
$arr[2,8] = 0.02807111181859822
$arr[2,9] = 0.02807111181859822
$arr[2,10] = 3
$arr[2,11] = 1
$arr[2,12] = 95.23175666666667
$arr[2,13] = 285.69527
$arr[2,14] = 0.2155966197102082
$arr[2,15] = 0.2155966197102082
$arr[2,16] = 13.70346885730667
$arr[2,17] = 123.33121971576
$arr[2,18] = 0.006052036819597053
$arr[2,19] = 0.006052036819597054
$arr[3,0] = "ECs"
$arr[3,1] = "Slit2"
$arr[3,2] = "App"
$arr[3,3] = "sCs"
$arr[3,4] = 2
$arr[3,5] = 0.6666666666666666
$arr[3,6] = 0.143896
$arr[3,7] = 0.431688
$arr[3,8] = 0.02807111181859822
$arr[3,9] = 0.02807111181859822
$arr[3,10] = 3
$arr[3,11] = 1
$arr[3,12] = 50.778675
$arr[3,13] = 152.336025
$arr[3,14] = 0.1149586132458188
$arr[3,15] = 0.1149586132458188
$arr[3,16] = 7.3068482178
$arr[3,17] = 65.7616339602
$arr[3,18] = 0.003227016086934366
$arr[3,19] = 0.003227016086934366
$arr[4,0] = "FAPs"
$arr[4,1] = "Slit2"
$arr[4,2] = "App"
$arr[4,3] = "ECs"
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 4.277274333333334
$arr[4,7] = 12.831823
$arr[4,8] = 0.8344071140950421
$arr[4,9] = 0.8344071140950421
$arr[4,10] = 3
$arr[4,11] = 1
$arr[4,12] = 110.8604276666667
$arr[4,13] = 332.581283
$arr[4,14] = 0.2509786052589675
$arr[4,15] = 0.2509786052589675
$arr[4,16] = 474.1804618409899
$arr[4,17] = 4267.624156568909
$arr[4,18] = 0.2094183337137338
$arr[4,19] = 0.2094183337137338
$arr[5,0] = "FAPs"
$arr[5,1] = "Slit2"
$arr[5,2] = "App"
$arr[5,3] = "FAPs"
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 4.277274333333334
$arr[5,7] = 12.831823
$arr[5,8] = 0.8344071140950421
$arr[5,9] = 0.8344071140950421
$arr[5,10] = 3
$arr[5,11] = 1
$arr[5,12] = 184.841802
$arr[5,13] = 554.525406
$arr[5,14] = 0.4184661617850055
$arr[5,15] = 0.4184661617850055
$arr[5,16] = 790.619095421682
$arr[5,17] = 7115.571858795138
$arr[5,18] = 0.3491711424014554
$arr[5,19] = 0.3491711424014554
$arr[6,0] = "FAPs"
$arr[6,1] = "Slit2"
$arr[6,2] = "App"
$arr[6,3] = "M2"
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = 4.277274333333334
$arr[6,7] = 12.831823
$arr[6,8] = 0.8344071140950421
$arr[6,9] = 0.8344071140950421
$arr[6,10] = 3
$arr[6,11] = 1
$arr[6,12] = 95.23175666666667
$arr[6,13] = 285.69527
$arr[6,14] = 0.2155966197102082
$arr[6,15] = 0.2155966197102082
$arr[6,16] = 407.3323485085789
$arr[6,17] = 3665.99113657721
$arr[6,18] = 0.1798953532610411
$arr[6,19] = 0.1798953532610411
$arr[7,0] = "FAPs"
$arr[7,1] = "Slit2"
$arr[7,2] = "App"
$arr[7,3] = "sCs"
$arr[7,4] = 3
$arr[7,5] = 1
$arr[7,6] = 4.277274333333334
$arr[7,7] = 12.831823
$arr[7,8] = 0.8344071140950421
$arr[7,9] = 0.8344071140950421
$arr[7,10] = 3
$arr[7,11] = 1
$arr[7,12] = 50.778675
$arr[7,13] = 152.336025
$arr[7,14] = 0.1149586132458188
$arr[7,15] = 0.1149586132458188
$arr[7,16] = 217.194323258175
$arr[7,17] = 1954.748909323575
$arr[7,18] = 0.09592228471881173
$arr[7,19] = 0.09592228471881174
$arr[8,0] = "sCs"
$arr[8,1] = "Slit2"
$arr[8,2] = "App"
$arr[8,3] = "ECs"
$arr[8,4] = 3
$arr[8,5] = 1
$arr[8,6] = 0.7049536666666666
$arr[8,7] = 2.114861
$arr[8,8] = 0.1375217740863597
$arr[8,9] = 0.1375217740863597
$arr[8,10] = 3
$arr[8,11] = 1
$arr[8,12] = 110.8604276666667
$arr[8,13] = 332.581283
$arr[8,14] = 0.2509786052589675
$arr[8,15] = 0.2509786052589675
$arr[8,16] = 78.15146497185144
$arr[8,17] = 703.3631847466629
$arr[8,18] = 0.03451502305293338
$arr[8,19] = 0.03451502305293339
$arr[9,0] = "sCs"
$arr[9,1] = "Slit2"
$arr[9,2] = "App"
$arr[9,3] = "FAPs"
$arr[9,4] = 3
$arr[9,5] = 1
$arr[9,6] = 0.7049536666666666
$arr[9,7] = 2.114861
$arr[9,8] = 0.1375217740863597
$arr[9,9] = 0.1375217740863597
$arr[9,10] = 3
$arr[9,11] = 1
$arr[9,12] = 184.841802
$arr[9,13] = 554.525406
$arr[9,14] = 0.4184661617850055
$arr[9,15] = 0.4184661617850055
$arr[9,16] = 130.304906073174
$arr[9,17] = 1172.744154658566
$arr[9,18] = 0.05754820896378357
$arr[9,19] = 0.05754820896378358
$arr[10,0] = "sCs"
$arr[10,1] = "Slit2"
$arr[10,2] = "App"
$arr[10,3] = "M2"
$arr[10,4] = 3
$arr[10,5] = 1
$arr[10,6] = 0.7049536666666666
$arr[10,7] = 2.114861
$arr[10,8] = 0.1375217740863597
$arr[10,9] = 0.1375217740863597
$arr[10,10] = 3
$arr[10,11] = 1
$arr[10,12] = 95.23175666666667
$arr[10,13] = 285.69527
$arr[10,14] = 0.2155966197102082
$arr[10,15] = 0.2155966197102082
$arr[10,16] = 67.13397604527444
$arr[10,17] = 604.20578440747
$arr[10,18] = 0.02964922962957006
$arr[10,19] = 0.02964922962957007
$arr[11,0] = "sCs"
$arr[11,1] = "Slit2"
$arr[11,2] = "App"
$arr[11,3] = "sCs"
$arr[11,4] = 3
$arr[11,5] = 1
$arr[11,6] = 0.7049536666666666
$arr[11,7] = 2.114861
$arr[11,8] = 0.1375217740863597
$arr[11,9] = 0.1375217740863597
$arr[11,10] = 3
$arr[11,11] = 1
$arr[11,12] = 50.778675
$arr[11,13] = 152.336025
$arr[11,14] = 0.1149586132458188
$arr[11,15] = 0.1149586132458188
$arr[11,16] = 35.796613129725
$arr[11,17] = 322.169518167525
$arr[11,18] = 0.01580931244007269
$arr[11,19] = 0.0158093124400727

$rng = $ws.Range("A2:T13")
$rng.Value = $arr
